$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 48, shifting rows 48-51 down to 49-52.
$ws.Range("A48:R48").Insert()

# The style used on column D (date) for the data rows is style index 2 (custom date format).
# Copy that style from the cell below (old row 48, now row 49) onto the new D48 cell.
$ws.Range("D49").Copy()
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new row 48 values.
$ws.Range("A48").Value = 10
$ws.Range("B48").Value = "Vega Modelo de Temuco"
$ws.Range("C48").Value = "La Araucanía"
$ws.Range("D48").Value = 44461
$ws.Range("E48").Value = 9
$ws.Range("F48").Value = 100112035
$ws.Range("G48").Value = "Bruselas (repollito)"
$ws.Range("H48").Value = "Sin especificar"
$ws.Range("I48").Value = "Primera"
$ws.Range("J48").Value = 40
$ws.Range("K48").Value = 25000
$ws.Range("L48").Value = 25000
$ws.Range("M48").Value = 25000
$ws.Range("N48").Value = "$/malla 10 kilos"
$ws.Range("O48").Value = "Provincia de Quillota"
$ws.Range("P48").Value = 2500
$ws.Range("Q48").Value = 10
$ws.Range("R48").Value = "Hortaliza"
